# laborator 12.12.2023 - Login part 2 + login pe frontend + add article
#
# Marks week 11 ("sapt 11", column M) attendance as present (TRUE) for the
# students that were present that week, mirroring what the other weekly
# attendance columns (C..L) already look like. The Q column ("Nota") is a
# shared formula that sums the weekly attendance, so it recalculates on its
# own once the M cells are populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where the student attended week 11 (column M)
$attendedRows = @(6, 10, 11, 17, 18, 21, 25, 28, 36, 37, 41, 42, 46)

foreach ($row in $attendedRows) {
    $ws.Cells.Item($row, 13).Value = $true   # column M = 13 ("sapt 11")
}

# Restore the view/selection state recorded for the sheet after the edit:
# scrolled so row 7 is at the top, with M11 the active selected cell.
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("M11").Select()
